# Insert a new record (row) right before the current row 331 of the
# "Hortaliza, Mercado Mayorista Lo Valledor de Santiago - Zapallo italiano"
# sheet. This shifts the existing rows 331..419 down to 332..420 and grows
# the used range from A1:R419 to A1:R420 (weekly update commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 331 (and everything below it) down by one row.
$ws.Rows.Item(331).Insert()

# Populate the newly inserted row 331 with the new data point.
$ws.Cells.Item(331,1).Value()  = 6
$ws.Cells.Item(331,2).Value()  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(331,3).Value()  = "Metropolitana"
$ws.Cells.Item(331,4).Value()  = 44754
$ws.Cells.Item(331,5).Value()  = 13
$ws.Cells.Item(331,6).Value()  = 100112032
$ws.Cells.Item(331,7).Value()  = "Zapallo italiano"
$ws.Cells.Item(331,8).Value()  = "Sin especificar"
$ws.Cells.Item(331,9).Value()  = "Primera"
$ws.Cells.Item(331,10).Value() = 200
$ws.Cells.Item(331,11).Value() = 7000
$ws.Cells.Item(331,12).Value() = 8000
$ws.Cells.Item(331,13).Value() = 7400
$ws.Cells.Item(331,14).Value() = "`$/caja 50 unidades"
$ws.Cells.Item(331,15).Value() = "Región de Arica y Parinacota"
$ws.Cells.Item(331,16).Value() = 148
$ws.Cells.Item(331,17).Value() = 50
$ws.Cells.Item(331,18).Value() = "Hortaliza"
